$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new requirement row ("1.2" - Evitar remoção de filial) -------
# The new row is inserted as row 6 (right after the existing "1.1" row),
# pushing every following row (old 6..68) down by one (new 7..69).
$ws.Rows.Item(6).EntireRow.Insert()

# --- Fill in the new row's content -----------------------------------------
$ws.Range("A6").Value = "RFN1.2"

# B6 ("1.2") looks like a number to Excel's parser, but in the source file
# it must be stored as text (shared string), like the sibling "1.1" cell
# (B5). Force the cell to text *before* assigning the value so the string
# is not auto-converted to a numeric 1.2.
$ws.Range("B6").NumberFormat = "@"
$ws.Range("B6").Value = "1.2"

# Copy the row-5 formatting (same visual "sub-item" row style used for the
# other "1.x" rows) onto the newly inserted row 6. This only touches
# formatting, so the text values assigned above are preserved.
$ws.Range("A5:I5").Copy()
$ws.Range("A6:I6").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("C6").Value = "Evitar remoção de filial."
$ws.Range("D6").Value = "Apresentação de mensagem de erro indicando “Existem alunos e colaboradores associados à filial.”."
$ws.Range("E6").Value = "Remover uma filial que contenha alunos ou colaboradores associados. "
$ws.Range("F6").Value = "Equipe de TCC"
$ws.Range("G6").Value = "Cliente"
$ws.Range("H6").Value = 43738
$ws.Range("I6").Value = "Sim"

# --- Update the view to match the saved cursor/scroll position -------------
$ws.Application.ActiveWindow.ScrollRow = 1
$sel = $ws.Range("D46")
$sel.Select()
